$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Opportunities
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Opportunities")
$ws1.Activate()

# Re-sort the holdings table (rows 5:11) back into ascending order by
# Security Code (column B) -- this undoes the earlier re-ordering.
$ws1.Range("B5:O11").Sort($ws1.Range("B5:B11"), 1)

# Rebuild the shared formulas for the calculated columns so they span
# the freshly sorted range again.
$ws1.Range("G5:G11").Formula = "=E5-H5"
$ws1.Range("N5:N11").Formula = "=M5/E5"

# Roll back the small data drift on the valuation inputs for each
# ticker (price / non-op assets / excess return / FCF value /
# realizable VPS / LFY dividend).
$ws1.Range("E5").Value2 = 20.8500003814697265625
$ws1.Range("I5").Value2 = -0.01676313671509410685
$ws1.Range("M5").Value2 = 0.17568614818334662986

$ws1.Range("E6").Value2 = 1.82000005245208740234
$ws1.Range("H6").Value2 = 0.80101504335149775748
$ws1.Range("I6").Value2 = -0.16284013534803840062
$ws1.Range("K6").Value2 = 1.5970421402726409088
$ws1.Range("L6").Value2 = 0.9934935520797943953
$ws1.Range("M6").Value2 = 0.05285658156263032398

$ws1.Range("E8").Value2 = 6.84999990463256835938
$ws1.Range("H8").Value2 = 0.57730142848580079384
$ws1.Range("I8").Value2 = -0.11019497191780214251
$ws1.Range("K8").Value2 = 6.60969463057438133546
$ws1.Range("L8").Value2 = 2.11397005095674339614
$ws1.Range("M8").Value2 = 0.02100100680006224824

$ws1.Range("E9").Value2 = 3.56999993324279785156
$ws1.Range("I9").Value2 = -0.03412413466086086644
$ws1.Range("M9").Value2 = 0.05820001292746511401

$ws1.Range("E10").Value2 = 0.60000002384185791016
$ws1.Range("I10").Value2 = -0.02242666211152535838

$ws1.Range("E11").Value2 = 1.89999997615814208984
$ws1.Range("I11").Value2 = 0.59498285502339665065
$ws1.Range("M11").Value2 = 0.0231580836679025695

# Update the sheet's remembered selection.
$ws1.Range("O10").Select()

# ---------------------------------------------------------------------
# Sheet 2: Current_Holdings
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Current_Holdings")
$ws2.Activate()
$ws2.Range("B7").Select()

# ---------------------------------------------------------------------
# Sheet 3: Discount rates
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Discount rates")
$ws3.Activate()
$excel.ActiveWindow.Zoom = 100
$ws3.Range("F53").Select()

# Restore the originally active sheet/tab.
$ws1.Activate()
